$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.617.79'
$ws.Range("E2").Value = '  -3.17%  '
$ws.Range("D3").Value = '2.095.41'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").Value = "'342.22"
$ws.Range("E5").Value = '  -2.36%  '
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = "'0.5120"
$ws.Range("E7").Value = '  -2.43%  '
$ws.Range("D8").Value = "'0.4396"
$ws.Range("E8").Value = '  -2.50%  '
$ws.Range("D9").Value = "'53.19"
$ws.Range("E9").Value = '  -2.20%  '
$ws.Range("D10").Value = "'0.09118"
$ws.Range("E10").Value = '  +1.24%  '
$ws.Range("D11").Value = "'1.167"
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("D12").Value = "'24.62"
$ws.Range("E12").Value = '  +0.75%  '
$ws.Range("D13").Value = '2.100.24'
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("D14").Value = "'6.748"
$ws.Range("D15").Value = "'8.187"
$ws.Range("E15").Value = '  +1.99%  '
$ws.Range("D16").Value = "'99.79"
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = "'0.00001146"
$ws.Range("E17").Value = '  -2.10%  '
$ws.Range("B18").Value = 'BinanceUSD'
$ws.Range("C18").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D18").Value = "'1.008"
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("D19").Value = "'21.01"
$ws.Range("E19").Value = '  +8.36%  '
$ws.Range("D20").Value = "'0.06645"
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("D21").Value = "'1.007"
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").Value = "'6.172"
$ws.Range("E22").Value = '  -1.91%  '
$ws.Range("D23").Value = '29.651.34'
$ws.Range("E23").Value = '  -3.28%  '
$ws.Range("D24").Value = "'12.57"
$ws.Range("E24").Value = '  -1.95%  '
$ws.Range("D25").Value = "'2.304"
$ws.Range("E25").Value = '  -3.47%  '
$ws.Range("D26").Value = '2.341.60'
$ws.Range("E26").Value = '  -1.13%  '
$ws.Range("D27").Value = "'21.81"
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("D28").Value = "'162.33"
$ws.Range("E28").Value = '  -1.67%  '
$ws.Range("D29").Value = "'2.517"
$ws.Range("E29").Value = '  -0.66%  '
$ws.Range("D30").Value = "'132.32"
$ws.Range("E30").Value = '  -2.78%  '
$ws.Range("D31").Value = "'1.127"
$ws.Range("E31").Value = '  -5.37%  '
$ws.Range("D32").Value = "'0.1044"
$ws.Range("E32").Value = '  -3.06%  '
$ws.Range("D33").Value = "'1.637"
$ws.Range("E33").Value = '  -0.89%  '
$ws.Range("D34").Value = "'6.139"
$ws.Range("E34").Value = '  -3.55%  '
$ws.Range("D35").Value = "'3.967"
$ws.Range("E35").Value = '  -1.40%  '
$ws.Range("D36").Value = "'6.044"
$ws.Range("E36").Value = '  +2.06%  '
$ws.Range("D37").Value = "'10.26"
$ws.Range("D38").Value = "'0.02563"
$ws.Range("E38").Value = '  -3.06%  '
$ws.Range("D39").Value = "'0.06656"
$ws.Range("E39").Value = '  -2.68%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = "'0.2224"
$ws.Range("E40").Value = '  -3.67%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = "'12.33"
$ws.Range("E41").Value = '  -1.50%  '
$ws.Range("D42").Value = "'0.6831"
$ws.Range("E42").Value = '  -0.62%  '
$ws.Range("D43").Value = "'1.285"
$ws.Range("E43").Value = '  +1.03%  '
$ws.Range("D44").Value = "'0.6637"
$ws.Range("E44").Value = '  +2.96%  '
$ws.Range("D45").Value = "'14.07"
$ws.Range("E45").Value = '  -4.52%  '
$ws.Range("D46").Value = "'2.287"
$ws.Range("E46").Value = '  -1.98%  '
$ws.Range("D47").Value = "'3.607"
$ws.Range("E47").Value = '  -3.96%  '
$ws.Range("D48").Value = "'1.218"
$ws.Range("E48").Value = '  -2.64%  '
$ws.Range("D49").Value = "'81.82"
$ws.Range("E49").Value = '  -0.73%  '
$ws.Range("D50").Value = "'0.00000000331"
$ws.Range("E50").Value = '  -7.49%  '
$ws.Range("D51").Value = "'1.160"
$ws.Range("E51").Value = '  -2.52%  '
